# Add two new side-by-side tables in rows 17-31 (titled "2020 United States
# presidential election" in A17:E17 and "COVID-19 pandemic" in G17:K17),
# mirroring the layout of the existing M17:Q31 table, and update the
# existing table's "titles" column (P19:P31) with revised figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Row 17: section headers, reusing the same text already used at the top
# of the sheet (A1 = "2020 United States presidential election",
# E1 = "COVID-19 pandemic").
# ---------------------------------------------------------------------
$ws.Range("A17").Value = $ws.Range("A1").Value
$ws.Range("G17").Value = $ws.Range("E1").Value

# Match the formatting of the existing M17:Q17 header band.
$ws.Range("M17:Q17").Copy()
$ws.Range("A17:E17").PasteSpecial($xlPasteFormats)
$ws.Range("M17:Q17").Copy()
$ws.Range("G17:K17").PasteSpecial($xlPasteFormats)

$ws.Range("A17:E17").Merge()
$ws.Range("G17:K17").Merge()

# ---------------------------------------------------------------------
# Row 18: column sub-headers.
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "nodes"
$ws.Range("C18").Value = "edges"
$ws.Range("D18").Value = "titles"
$ws.Range("E18").Value = "categories"
$ws.Range("H18").Value = "nodes"
$ws.Range("I18").Value = "edges"
$ws.Range("J18").Value = "titles"
$ws.Range("K18").Value = "categories"

$ws.Range("M18:Q18").Copy()
$ws.Range("A18:E18").PasteSpecial($xlPasteFormats)
$ws.Range("M18:Q18").Copy()
$ws.Range("G18:K18").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Rows 19-31: data. Each entry is
#   date, nodes, edges, titles, categories (left table, A:E)
#   nodes, edges, titles, categories (right table, H:K -- date mirrors A)
#   revised "titles" figure for the pre-existing M:Q table (column P)
# ---------------------------------------------------------------------
$rows = @(
  @("Jan. 1st",  1220, 90465,  5236,  10773, 1778, 17890,  6616,  11586, 2330),
  @("Jan. 20th", 1226, 93924,  5269,  10880, 1853, 18982,  6952,  12098, 2394),
  @("Feb. 1st",  1229, 95755,  5280,  10946, 1911, 19876,  7191,  12576, 2424),
  @("Feb. 20th", 1245, 98824,  5378,  11094, 2040, 22456,  7809,  13289, 2445),
  @("Mar. 1st",  1250, 99800,  5400,  11148, 2120, 25485,  8185,  13760, 2464),
  @("Apr. 1st",  1276, 102608, 5503,  11338, 2759, 208020, 11296, 18008, 2534),
  @("Apr. 20th", 1296, 105270, 5636,  11527, 3006, 276751, 12334, 20470, 2554),
  @("May 1st",   1303, 106657, 5670,  11583, 3127, 314667, 12789, 21503, 2565),
  @("May 20th",  1311, 108050, 5690,  11663, 3278, 367660, 12884, 23024, 2584),
  @("Jun. 1st",  1314, 108421, 5715,  11751, 3351, 382105, 13163, 23767, 2722),
  @("Jun. 20th", 1322, 108704, 5843,  11865, 3431, 405540, 13520, 24572, 3530),
  @("Jul. 1st",  1331, 110250, 5884,  11953, 3471, 429762, 13688, 25120, 3688),
  @("Jul. 20th", 1342, 112344, 5935,  12047, 3554, 456680, 14032, 25946, 3812)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 19 + $i
    $row = $rows[$i]

    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
    $ws.Range("E" + $r).Value = $row[4]

    $ws.Range("G" + $r).Value = $row[0]
    $ws.Range("H" + $r).Value = $row[5]
    $ws.Range("I" + $r).Value = $row[6]
    $ws.Range("J" + $r).Value = $row[7]
    $ws.Range("K" + $r).Value = $row[8]

    $srcRange = "M" + $r + ":Q" + $r
    $ws.Range($srcRange).Copy()
    $ws.Range("A" + $r + ":E" + $r).PasteSpecial($xlPasteFormats)
    $ws.Range($srcRange).Copy()
    $ws.Range("G" + $r + ":K" + $r).PasteSpecial($xlPasteFormats)

    $ws.Range("P" + $r).Value = $row[9]
}

# ---------------------------------------------------------------------
# Restore the view/selection state recorded in the workbook.
# ---------------------------------------------------------------------
$null = $ws.Range("A10").Select()
$null = $ws.Range("I32").Select()

Write-Host "done"
